$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row contents (A/B columns) for rows 6-12 ---
$ws.Range("A6").Value = "bitcoin beginners  "
$ws.Range("B6").Value = "com.hamxa.shaynachim"

$ws.Range("A7").Value = "passive income"
$ws.Range("B7").Value = "passive.income.nadi.myfirstdrawermenuproject"

$ws.Range("A8").Value = "Best bitcoin"
$ws.Range("B8").Value = "com.hamxa.shaynachim"

$ws.Range("A9").Value = "bitcoin"
$ws.Range("B9").Value = "com.hamxa.shaynachim"

$ws.Range("A10").Value = "powerful quotes"
$ws.Range("B10").Value = "com.sugar.powerfulquotes"

$ws.Range("A11").Value = "bitcoin"
$ws.Range("B11").Value = "com.hamxa.shaynachim"

$ws.Range("A12").Value = "blockchain"
$ws.Range("B12").Value = "block.chain.technology"

# --- Row height changes: row 6 grows to 24, rows 8 & 9 shrink back to default 12.8 ---
$ws.Rows.Item(6).RowHeight = 24
$ws.Rows.Item(8).RowHeight = 12.8
$ws.Rows.Item(9).RowHeight = 12.8

# --- Update the saved selection / active cell ---
$ws.Range("B22").Select() | Out-Null
